$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04912295720824748
$ws.Range("D2").Value = 0.9612645567762683

$ws.Range("C3").Value = -0.2037890831511524
$ws.Range("D3").Value = 0.8403916250389059

$ws.Range("C4").Value = -0.07903642456770558
$ws.Range("D4").Value = 0.9377181626776185

$ws.Range("C5").Value = 0.4926943900923948
$ws.Range("D5").Value = 0.6271099885974034

$ws.Range("C6").Value = -0.2690620020950505
$ws.Range("D6").Value = 0.7903892293056161

$ws.Range("C7").Value = -0.1404202391429826
$ws.Range("D7").Value = 0.889605699955456

$ws.Range("C8").Value = 1.088164056914269
$ws.Range("D8").Value = 0.2883004501516211

$ws.Range("C9").Value = 0.184950544150321
$ws.Range("D9").Value = 0.8549622477527479

$ws.Range("C10").Value = 0.8261477443632111
$ws.Range("D10").Value = 0.4175901721580413

$ws.Range("C11").Value = 0.7530079019237732
$ws.Range("D11").Value = 0.4594278865819785
